$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix skill-name / image-name correspondence (swap columns) ---
# Rows 12-15: swap F (articulation/dynamics) <-> G (dynamics/articulation)
foreach ($r in 12..15) {
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $tmp = $fCell.Value2
    $fCell.Value2 = $gCell.Value2
    $gCell.Value2 = $tmp
}

# Rows 16-21: swap H (stim_a/stim_d) <-> I (stim_d/stim_a)
foreach ($r in 16..21) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $iCell = $ws.Cells.Item($r, 9)   # column I
    $tmp = $hCell.Value2
    $hCell.Value2 = $iCell.Value2
    $iCell.Value2 = $tmp
}

# --- Remove the last two (now duplicated / unnecessary) trial rows ---
$ws.Rows("30:31").Delete()

# --- Keep the print area in sync with the shrunk data range ---
$ws.PageSetup.PrintArea = "`$A`$1:`$I`$29"

# --- Restore the selection to its saved cell ---
[void]$ws.Range("G8").Select()

# --- Page setup: fit-to-page flag on, with an 84% scale remembered ---
$ws.PageSetup.Zoom = 84
$ws.PageSetup.FitToPagesWide = 0
$ws.PageSetup.FitToPagesTall = 0
